# employees.xlsx import-script improvement
# - Portal "Active" column (R) used to store a raw boolean (TRUE/FALSE);
#   switch it to the human-readable "Active"/"Inactive" text the importer
#   now emits.
# - Normalize the "Mr. " title (trailing period+space) down to "Mr" so it
#   matches the other bare titles ("Dr") already in the sheet.
# - Leave the last selected cell where the author's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("employees")

$ws.Range("R2").Value = "Active"
$ws.Range("R3").Value = "Active"
$ws.Range("R4").Value = "Inactive"
$ws.Range("O5").Value = "Mr"
$ws.Range("R5").Value = "Active"

$ws.Range("H17").Select() | Out-Null
